$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.0498220640569395
$summary.Range("C2").Value = 0.0498220640569395
$summary.Range("D2").Value = 1
$summary.Range("E2").Value = 0.09491525423728814
$summary.Range("F2").Value = 0.2077151335311573
$summary.Range("G2").Value = 0.5768621236133122
$summary.Range("H2").Value = 0.6927166934189406
$summary.Range("I2").Value = 28
$summary.Range("J2").Value = 534
$summary.Range("K2").Value = 0
$summary.Range("L2").Value = 0

# --- Sheet "Classification Report" ---
$cr = $wb.Worksheets.Item("Classification Report")

# row 2 ("0")
$cr.Range("B2").Value = 0
$cr.Range("C2").Value = 0
$cr.Range("D2").Value = 0

# row 3 ("1")
$cr.Range("B3").Value = 0.0498220640569395
$cr.Range("C3").Value = 1
$cr.Range("D3").Value = 0.09491525423728814

# row 4 ("accuracy")
$cr.Range("B4").Value = 0.0498220640569395
$cr.Range("C4").Value = 0.0498220640569395
$cr.Range("D4").Value = 0.0498220640569395
$cr.Range("E4").Value = 0.0498220640569395

# row 5 ("macro avg")
$cr.Range("B5").Value = 0.02491103202846975
$cr.Range("C5").Value = 0.5
$cr.Range("D5").Value = 0.04745762711864407

# row 6 ("weighted avg")
$cr.Range("B6").Value = 0.002482238066893783
$cr.Range("C6").Value = 0.0498220640569395
$cr.Range("D6").Value = 0.004728873876590867

# --- Sheet "Confusion Matrix" ---
$cm = $wb.Worksheets.Item("Confusion Matrix")

# row 2 ("Actual 0")
$cm.Range("B2").Value = 0
$cm.Range("C2").Value = 534

# row 3 ("Actual 1")
$cm.Range("B3").Value = 0
$cm.Range("C3").Value = 28
